$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.525.64"
$ws.Range("E2").Value = "  +5.55%  "

$ws.Range("D3").Value = "1.707.37"
$ws.Range("E3").Value = "  +4.15%  "

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "'221.86"
$ws.Range("E5").Value = "  +2.99%  "

$ws.Range("D6").Value = "'0.536"
$ws.Range("E6").Value = "  +3.22%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").Value = "'29.93"
$ws.Range("E8").Value = "  +4.25%  "

$ws.Range("D9").Value = "'45.35"
$ws.Range("E9").Value = "  +3.40%  "

$ws.Range("D10").Value = "'0.270"
$ws.Range("E10").Value = "  +3.58%  "

$ws.Range("E11").Value = "  +5.73%  "

$ws.Range("D13").Value = "1.951.61"
$ws.Range("E13").Value = "  +4.08%  "

$ws.Range("D14").Value = "1.701.46"
$ws.Range("E14").Value = "  +3.75%  "

$ws.Range("D15").Value = "'10.31"
$ws.Range("E15").Value = "  +9.14%  "

$ws.Range("D16").Value = "'0.613"
$ws.Range("E16").Value = "  +3.85%  "

$ws.Range("E17").Value = "  +8.18%  "

$ws.Range("D18").Value = "31.447.91"
$ws.Range("E18").Value = "  +5.24%  "

$ws.Range("D19").Value = "'67.21"
$ws.Range("E19").Value = "  +4.02%  "

$ws.Range("D20").Value = "'251.30"
$ws.Range("E20").Value = "  +4.48%  "

$ws.Range("D21").Value = "0.0₃0727"
$ws.Range("E21").Value = "  +3.44%  "

$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").Value = "'10.22"
$ws.Range("E23").Value = "  +3.30%  "

$ws.Range("E24").Value = "  +3.31%  "

$ws.Range("E25").Value = "  -1.57%  "

$ws.Range("D26").Value = "'159.45"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").Value = "'16.05"
$ws.Range("E27").Value = "  +3.23%  "

$ws.Range("E28").Value = "  +3.14%  "

$ws.Range("D29").Value = "'6.82"
$ws.Range("E29").Value = "  +2.96%  "

$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("E31").Value = "  +11.76%  "

$ws.Range("E32").Value = "  +2.22%  "

$ws.Range("E33").Value = "  +3.88%  "

$ws.Range("E34").Value = "  +6.92%  "

$ws.Range("D35").Value = "1.511.51"
$ws.Range("E35").Value = "  +6.06%  "

$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("E37").Value = "  +2.01%  "

$ws.Range("D38").Value = "'83.67"
$ws.Range("E38").Value = "  +9.24%  "

$ws.Range("E39").Value = "  +8.83%  "

$ws.Range("E40").Value = "  +4.09%  "

$ws.Range("D41").Value = "'2.69"
$ws.Range("E41").Value = "  -1.45%  "

$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").Value = "'2.05"
$ws.Range("E43").Value = "  +4.04%  "

$ws.Range("E44").Value = "  +2.35%  "

$ws.Range("D45").Value = "'0.0503"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").Value = "'52.65"
$ws.Range("E48").Value = "  +7.13%  "

$ws.Range("E49").Value = "  +3.76%  "

$ws.Range("E50").Value = "  +3.55%  "

$ws.Range("D51").Value = "0.0₆0120"
$ws.Range("E51").Value = "  +10.53%  "

